$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.5945692883895131
$ws1.Range("C2").Value = 0.5524402907580478
$ws1.Range("D2").Value = 0.9962546816479401
$ws1.Range("E2").Value = 0.7107548430193721
$ws1.Range("F2").Value = 0.8583414004517587
$ws1.Range("G2").Value = 0.9663941871026339
$ws1.Range("H2").Value = 0.7628596277125504
$ws1.Range("I2").Value = 532
$ws1.Range("J2").Value = 431
$ws1.Range("K2").Value = 103
$ws1.Range("L2").Value = 2

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.9809523809523809
$ws2.Range("C2").Value = 0.1928838951310861
$ws2.Range("D2").Value = 0.3223787167449139

$ws2.Range("B3").Value = 0.5524402907580478
$ws2.Range("C3").Value = 0.9962546816479401
$ws2.Range("D3").Value = 0.7107548430193721

$ws2.Range("B4").Value = 0.5945692883895131
$ws2.Range("C4").Value = 0.5945692883895131
$ws2.Range("D4").Value = 0.5945692883895131
$ws2.Range("E4").Value = 0.5945692883895131

$ws2.Range("B5").Value = 0.7666963358552144
$ws2.Range("C5").Value = 0.5945692883895131
$ws2.Range("D5").Value = 0.5165667798821429

$ws2.Range("B6").Value = 0.7666963358552144
$ws2.Range("C6").Value = 0.5945692883895131
$ws2.Range("D6").Value = 0.5165667798821431

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 103
$ws3.Range("C2").Value = 431
$ws3.Range("B3").Value = 2
$ws3.Range("C3").Value = 532
